$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E25").Value = "['Normal', 'SoftwareFault']"

# Row 38
$ws.Range("D38").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['HardwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault']"

# Row 56
$ws.Range("D56").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['HardwareFault']"

# Row 58
$ws.Range("D58").Value = "[0, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['ParamViolation']"

# Row 69
$ws.Range("D69").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E69").Value = "['Normal', 'SurroundingEnvironment']"

# Row 74
$ws.Range("D74").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E74").Value = "['Normal']"

# Row 84
$ws.Range("D84").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E84").Value = "['Normal', 'SurroundingEnvironment']"

# Row 91
$ws.Range("D91").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E91").Value = "['Normal', 'HardwareFault']"

# Row 107
$ws.Range("D107").Value = "[1, 0, 0, 0, 0, 1, 0]"
$ws.Range("E107").Value = "['Normal', 'CommunicationIssue']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E113").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"
